$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (14 and 15) by copying the formatting/pattern of the
# existing "pair" row 12 (style pattern: B/F = style 2, C/D/E = style 1).
# Doing the insert twice within the original used range (row <= 14) makes the
# engine carry over the per-cell styles (s="2"/s="1") to the newly created rows.
$ws.Rows(12).Copy()
$ws.Rows(14).Insert()
$ws.Rows(12).Copy()
$ws.Rows(14).Insert()

# Remove the row that got pushed down to 16 (the old, now-duplicated, blank
# row 14 from the original sheet); it is not part of the target layout.
$ws.Rows(16).Delete()

# Row 14: new single-sided file f6_1 (no related_files.file value)
$ws.Range("B14").Value = "test_lab:f6_1, test_lab:alt_f6_1"
$ws.Range("C14").Value = "fastq"
$ws.Range("D14").Value = "1"
$ws.Range("E14").Value = "paired with"
$ws.Range("F14").Value = $null
$ws.Range("F14").NumberFormat = "@"

# Row 15: new file f6_2, pointing back at f6_1
$ws.Range("B15").Value = "test_lab:f6_2, test_lab:alt_f6_2"
$ws.Range("C15").Value = "fastq"
$ws.Range("D15").Value = "2"
$ws.Range("E15").Value = "paired with"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "test_lab:f6_1"

# Update the selection to match the saved view state.
$ws.Range("B26").Select()
